$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.477963924407959
$ws.Range("B1").Value = 3.239269256591797
$ws.Range("C1").Value = 2.845316886901855
$ws.Range("D1").Value = 3.09270191192627
$ws.Range("E1").Value = 2.373469352722168
